# Auto-generated edit script applying numeric corrections to Titan_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1020.75
$ws.Range("I18").Value = 333.4
$ws.Range("K18").Value = 333.4
$ws.Range("M18").Value = -49.39999999999998
$ws.Range("H74").Value = 3288
$ws.Range("I74").Value = 2433.3333
$ws.Range("K74").Value = 2433.3333
$ws.Range("M74").Value = -1497.3333
$ws.Range("H77").Value = 3288
$ws.Range("I77").Value = 2433.3333
$ws.Range("K77").Value = 12166.6665
$ws.Range("M77").Value = -7486.666499999999
$ws.Range("H118").Value = 1320
$ws.Range("J118").Value = 1600
$ws.Range("L118").Value = 4800
$ws.Range("N118").Value = -8114
$ws.Range("H127").Value = 1177.125
$ws.Range("I127").Value = 787.5
$ws.Range("J127").Value = 1307
$ws.Range("K127").Value = 2362.5
$ws.Range("L127").Value = 3921
$ws.Range("M127").Value = 2597.5
$ws.Range("N127").Value = -13841
$ws.Range("H129").Value = 712.1667
$ws.Range("I129").Value = 293.25
$ws.Range("J129").Value = 1550
$ws.Range("K129").Value = 879.75
$ws.Range("L129").Value = 4650
$ws.Range("M129").Value = 4120.25
$ws.Range("N129").Value = -14650
$ws.Range("H132").Value = 15979.913
$ws.Range("I132").Value = 16539.605
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 49618.815
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -47088.815
$ws.Range("N132").Value = -16060.0001
$ws.Range("H135").Value = 1076.0952
$ws.Range("I135").Value = 1171.7778
$ws.Range("K135").Value = 10546.0002
$ws.Range("M135").Value = -8011.0002
$ws.Range("H138").Value = 11336618
$ws.Range("I138").Value = 3033496.2
$ws.Range("J138").Value = 17860498
$ws.Range("K138").Value = 9100488.600000001
$ws.Range("L138").Value = 53581494
$ws.Range("M138").Value = -9095348.600000001
$ws.Range("N138").Value = -53591774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19818.852
$ws.Range("I32").Value = 4546.4697
$ws.Range("J32").Value = 145816
$ws.Range("K32").Value = 4546.4697
$ws.Range("L32").Value = 145816
$ws.Range("M32").Value = -4259.4697
$ws.Range("N32").Value = -146390
$ws.Range("H45").Value = 1064.2858
$ws.Range("I45").Value = 1075
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1075
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -698
$ws.Range("N45").Value = -1754
$ws.Range("H61").Value = 2906.65
$ws.Range("I61").Value = 2486.0667
$ws.Range("J61").Value = 4168.4
$ws.Range("K61").Value = 2486.0667
$ws.Range("L61").Value = 4168.4
$ws.Range("M61").Value = -2274.0667
$ws.Range("N61").Value = -4592.4
$ws.Range("H122").Value = 16303.125
$ws.Range("I122").Value = 18203.572
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 54610.716
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -52160.716
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 3439.8276
$ws.Range("I132").Value = 2908.6956
$ws.Range("J132").Value = 5475.8335
$ws.Range("K132").Value = 8726.086800000001
$ws.Range("L132").Value = 16427.5005
$ws.Range("M132").Value = -6196.086800000001
$ws.Range("N132").Value = -21487.5005
$ws.Range("H136").Value = 2906.65
$ws.Range("I136").Value = 2486.0667
$ws.Range("J136").Value = 4168.4
$ws.Range("K136").Value = 7458.2001
$ws.Range("L136").Value = 12505.2
$ws.Range("M136").Value = -4908.2001
$ws.Range("N136").Value = -17605.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5004799
$ws.Range("I7").Value = 6398.6665
$ws.Range("J7").Value = 20000000
$ws.Range("K7").Value = 6398.6665
$ws.Range("L7").Value = 20000000
$ws.Range("M7").Value = -6285.6665
$ws.Range("N7").Value = -20000226
$ws.Range("H82").Value = 13494.637
$ws.Range("I82").Value = 6271.222
$ws.Range("K82").Value = 6271.222
$ws.Range("M82").Value = -5888.222
$ws.Range("H85").Value = 13494.637
$ws.Range("I85").Value = 6271.222
$ws.Range("K85").Value = 6271.222
$ws.Range("M85").Value = -4945.222
$ws.Range("H94").Value = 839.72
$ws.Range("I94").Value = 673.8333
$ws.Range("K94").Value = 673.8333
$ws.Range("M94").Value = -222.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 50002750
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 52634420
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 52634420
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -52634644
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H122").Value = 1840.7931
$ws.Range("I122").Value = 1127.4
$ws.Range("J122").Value = 3426.111
$ws.Range("K122").Value = 3382.2
$ws.Range("L122").Value = 10278.333
$ws.Range("M122").Value = -932.2000000000003
$ws.Range("N122").Value = -15178.333
$ws.Range("H124").Value = 30000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 30000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 30000
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -34910
$ws.Range("H132").Value = 1654.9762
$ws.Range("I132").Value = 1129.6945
$ws.Range("J132").Value = 4806.6665
$ws.Range("K132").Value = 3389.0835
$ws.Range("L132").Value = 14419.9995
$ws.Range("M132").Value = -859.0835000000002
$ws.Range("N132").Value = -19479.9995
$ws.Range("H134").Value = 3286.8
$ws.Range("I134").Value = 1851.4445
$ws.Range("J134").Value = 6977.7144
$ws.Range("K134").Value = 5554.333500000001
$ws.Range("L134").Value = 20933.1432
$ws.Range("M134").Value = -3019.333500000001
$ws.Range("N134").Value = -26003.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H22").Value = 51508
$ws.Range("I22").Value = 51508
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 51508
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -50979
$ws.Range("N22").ClearContents()
$ws.Range("H123").Value = 11539.546
$ws.Range("J123").Value = 11539.546
$ws.Range("L123").Value = 11539.546
$ws.Range("N123").Value = -16439.546
$ws.Range("H132").Value = 3909.76
$ws.Range("I132").Value = 3335.9
$ws.Range("J132").Value = 6205.2
$ws.Range("K132").Value = 10007.7
$ws.Range("L132").Value = 18615.6
$ws.Range("M132").Value = -7477.700000000001
$ws.Range("N132").Value = -23675.6
$ws.Range("H136").Value = 21444.154
$ws.Range("J136").Value = 21379.334
$ws.Range("L136").Value = 64138.00199999999
$ws.Range("N136").Value = -69238.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2500
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H132").Value = 2576.111
$ws.Range("I132").Value = 1694.0834
$ws.Range("K132").Value = 5082.2502
$ws.Range("M132").Value = -2552.2502
$ws.Range("H136").Value = 3112.8057
$ws.Range("I136").Value = 1425.2188
$ws.Range("J136").Value = 16613.5
$ws.Range("K136").Value = 4275.6564
$ws.Range("L136").Value = 49840.5
$ws.Range("M136").Value = -1725.6564
$ws.Range("N136").Value = -54940.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 6042.2856
$ws.Range("J40").Value = 6666
$ws.Range("L40").Value = 6666
$ws.Range("N40").Value = -6964
$ws.Range("H122").Value = 1925.9375
$ws.Range("I122").Value = 1671.2
$ws.Range("K122").Value = 5013.6
$ws.Range("M122").Value = -2563.6
$ws.Range("H132").Value = 1946.3715
$ws.Range("I132").Value = 1789.5763
$ws.Range("J132").Value = 2787.3635
$ws.Range("K132").Value = 5368.7289
$ws.Range("L132").Value = 8362.0905
$ws.Range("M132").Value = -2838.7289
$ws.Range("N132").Value = -13422.0905
